# ALZ Policy Assignments v2.xlsx
# Insert a new "Microsoft Cloud Security Benchmark v2" initiative assignment
# row into the "ALZ Default Policy Assignments" sheet, as the new row 7
# (pushing all subsequent rows down by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALZ Default Policy Assignments")

# --- Insert a new blank row at position 7 (existing rows 7..74 shift to 8..75) ---
$ws.Rows.Item(7).Insert()

# --- Populate the new row 7 ---
# Column order: A Assignment Scope (MG), B Assignment Name, C Definition Name,
# D Type, E Custom/Builtin, F Description, G Effect, H GitHub Assignment File,
# I AzAdvertizer Link, J Release
$ws.Cells.Item(7, 1).Value2 = "Intermediate Root"
$ws.Cells.Item(7, 2).Value2 = "Microsoft Cloud Security Benchmark v2"
$ws.Cells.Item(7, 3).Value2 = "[Preview]: Microsoft cloud security benchmark v2"
$ws.Cells.Item(7, 4).Value2 = "Initiative"
$ws.Cells.Item(7, 5).Value2 = "Built-in"
$ws.Cells.Item(7, 8).Value2 = "DINE-ASB2PolicyAssignment.json"
$ws.Cells.Item(7, 6).Value2 = "The Microsoft cloud security benchmark initiative represents the policies and controls implementing security recommendations defined in Microsoft cloud security benchmark, see https://aka.ms/azsecbm. This also serves as the Microsoft Defender for Cloud default policy initiative. You can directly assign this initiative, or manage its policies and compliance results within Microsoft Defender for Cloud."
$ws.Cells.Item(7, 7).Value2 = "Audit, AuditIfNotExists, Disabled"
$ws.Cells.Item(7, 9).Value2 = "https://www.azadvertizer.net/azpolicyinitiativesadvertizer/e3ec7e09-768c-4b64-882c-fcada3772047.html"
$ws.Cells.Item(7, 10).Value2 = 45992

# Row 7 needs a taller row height to fit the new (longer) wrapped description.
$ws.Rows.Item(7).RowHeight = 144

# --- Fix up the hyperlinks whose target cells shifted down by one row ---
# (row 73 -> 74, row 74 -> 75, row 17 -> 18)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("I74"), "https://www.azadvertizer.net/azpolicyinitiativesadvertizer/281d9e47-d14d-4f05-b8eb-18f2c4a034ff.html")
$ws.Hyperlinks.Add($ws.Range("I75"), "https://www.azadvertizer.net/azpolicyinitiativesadvertizer/281d9e47-d14d-4f05-b8eb-18f2c4a034ff.html")
$ws.Hyperlinks.Add($ws.Range("I18"), "https://www.azadvertizer.net/azpolicyadvertizer/98903777-a9f6-47f5-90a9-acaf62ab01a8.html")

# --- Restore the view state (scroll position / active cell) ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("K7").Select()
